$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 152
$ws.Range("I2").Value = 450
$ws.Range("J2").Value = 1771
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 467
$ws.Range("M2").Value = 32
$ws.Range("N2").Value = 297
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 26
$ws.Range("S2").Value = 186
$ws.Range("T2").Value = 321
$ws.Range("U2").Value = 22
$ws.Range("V2").Value = 2669
$ws.Range("X2").Value = 2703
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 43
$ws.Range("AA2").Value = 14
